$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3022.9167
$ws.Range("I100").Value = 3466.3333
$ws.Range("J100").Value = 2875.111
$ws.Range("K100").Value = 3466.3333
$ws.Range("L100").Value = 2875.111
$ws.Range("M100").Value = -2925.3333
$ws.Range("N100").Value = -3957.111
$ws.Range("H116").Value = 8185.5
$ws.Range("I116").Value = 9875.299999999999
$ws.Range("J116").Value = 3961
$ws.Range("K116").Value = 9875.299999999999
$ws.Range("L116").Value = 3961
$ws.Range("M116").Value = -6433.299999999999
$ws.Range("N116").Value = -10845
$ws.Range("H132").Value = 10162756
$ws.Range("I132").Value = 10449092
$ws.Range("K132").Value = 31347276
$ws.Range("M132").Value = -31344746
$ws.Range("H137").Value = 13687.5
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 13687.5
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 41062.5
$ws.Range("M137").Value = ""
$ws.Range("N137").Value = -46162.5
$ws.Range("H138").Value = 2570.89
$ws.Range("J138").Value = 3125.743
$ws.Range("L138").Value = 9377.228999999999
$ws.Range("N138").Value = -19657.229

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 9768
$ws.Range("J27").Value = 9768
$ws.Range("L27").Value = 9768
$ws.Range("N27").Value = -10136
$ws.Range("H74").Value = 12441.556
$ws.Range("I74").Value = 1425.1428
$ws.Range("J74").Value = 50999
$ws.Range("K74").Value = 1425.1428
$ws.Range("L74").Value = 50999
$ws.Range("M74").Value = -551.1428000000001
$ws.Range("N74").Value = -52747
$ws.Range("H77").Value = 12441.556
$ws.Range("I77").Value = 1425.1428
$ws.Range("J77").Value = 50999
$ws.Range("K77").Value = 7125.714
$ws.Range("L77").Value = 254995
$ws.Range("M77").Value = -2757.714
$ws.Range("N77").Value = -263731
$ws.Range("H97").Value = 1187
$ws.Range("I97").Value = 1203.44
$ws.Range("K97").Value = 1203.44
$ws.Range("M97").Value = -707.4400000000001
$ws.Range("H110").Value = 8669.781000000001
$ws.Range("I110").Value = 10891.263
$ws.Range("K110").Value = 10891.263
$ws.Range("M110").Value = -8846.263000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 16786
$ws.Range("J58").Value = 16786
$ws.Range("L58").Value = 16786
$ws.Range("N58").Value = -17374
$ws.Range("H86").Value = 8233.609
$ws.Range("I86").Value = 6068.5864
$ws.Range("J86").Value = 13465.75
$ws.Range("K86").Value = 6068.5864
$ws.Range("L86").Value = 13465.75
$ws.Range("M86").Value = -4945.5864
$ws.Range("N86").Value = -15711.75
$ws.Range("H89").Value = 8233.609
$ws.Range("I89").Value = 6068.5864
$ws.Range("J89").Value = 13465.75
$ws.Range("K89").Value = 30342.932
$ws.Range("L89").Value = 67328.75
$ws.Range("M89").Value = -24726.932
$ws.Range("N89").Value = -78560.75
$ws.Range("H107").Value = 3076.8462
$ws.Range("I107").Value = 3057
$ws.Range("K107").Value = 3057
$ws.Range("M107").Value = -1137
$ws.Range("H134").Value = 2362.5386
$ws.Range("I134").Value = 1693.6666
$ws.Range("J134").Value = 3867.5
$ws.Range("K134").Value = 5080.9998
$ws.Range("L134").Value = 11602.5
$ws.Range("M134").Value = -2545.9998
$ws.Range("N134").Value = -16672.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 256.8889
$ws.Range("J22").Value = 197
$ws.Range("L22").Value = 197
$ws.Range("N22").Value = -897
$ws.Range("H31").Value = 7287389
$ws.Range("J31").Value = 34098.445
$ws.Range("L31").Value = 34098.445
$ws.Range("N31").Value = -34688.445
$ws.Range("H34").Value = 7287389
$ws.Range("J34").Value = 34098.445
$ws.Range("L34").Value = 34098.445
$ws.Range("N34").Value = -34502.445
$ws.Range("H98").Value = 104995
$ws.Range("J98").Value = 104995
$ws.Range("L98").Value = 104995
$ws.Range("N98").Value = -109487
$ws.Range("H99").Value = 3939.3
$ws.Range("I99").Value = 3821.4443
$ws.Range("K99").Value = 3821.4443
$ws.Range("M99").Value = -2323.4443
$ws.Range("H105").Value = 2302
$ws.Range("I105").Value = 2286.5
$ws.Range("J105").Value = 2333
$ws.Range("K105").Value = 2286.5
$ws.Range("L105").Value = 2333
$ws.Range("M105").Value = -539.5
$ws.Range("N105").Value = -5827
$ws.Range("H107").Value = 1723.125
$ws.Range("I107").Value = 1880.8334
$ws.Range("J107").Value = 1250
$ws.Range("K107").Value = 1880.8334
$ws.Range("L107").Value = 1250
$ws.Range("M107").Value = 39.16660000000002
$ws.Range("N107").Value = -5090
$ws.Range("H122").Value = 2005.2222
$ws.Range("I122").Value = 2105.5
$ws.Range("J122").Value = 1804.6666
$ws.Range("K122").Value = 6316.5
$ws.Range("L122").Value = 5413.9998
$ws.Range("M122").Value = -3866.5
$ws.Range("N122").Value = -10313.9998
$ws.Range("H126").Value = 3939.3
$ws.Range("I126").Value = 3821.4443
$ws.Range("K126").Value = 11464.3329
$ws.Range("M126").Value = -8994.332900000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 11442.223
$ws.Range("I136").Value = 4015
$ws.Range("J136").Value = 13564.286
$ws.Range("K136").Value = 12045
$ws.Range("L136").Value = 40692.858
$ws.Range("M136").Value = -6945
$ws.Range("N136").Value = -50892.858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 31000
$ws.Range("J47").Value = 31000
$ws.Range("L47").Value = 31000
$ws.Range("N47").Value = -32136
$ws.Range("H55").Value = 10054.444
$ws.Range("J55").Value = 22000
$ws.Range("L55").Value = 22000
$ws.Range("N55").Value = -22654
$ws.Range("H70").Value = 6391.607
$ws.Range("I70").Value = 5983.778
$ws.Range("J70").Value = 7125.7
$ws.Range("K70").Value = 5983.778
$ws.Range("L70").Value = 7125.7
$ws.Range("M70").Value = -5713.778
$ws.Range("N70").Value = -7665.7
$ws.Range("H73").Value = 6391.607
$ws.Range("I73").Value = 5983.778
$ws.Range("J73").Value = 7125.7
$ws.Range("K73").Value = 5983.778
$ws.Range("L73").Value = 7125.7
$ws.Range("M73").Value = -5047.778
$ws.Range("N73").Value = -8997.700000000001
$ws.Range("H118").Value = 19982.334
$ws.Range("J118").Value = 19982.334
$ws.Range("L118").Value = 19982.334
$ws.Range("N118").Value = -23296.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 13926
$ws.Range("I7").Value = 14798.091
$ws.Range("J7").Value = 4333
$ws.Range("K7").Value = 14798.091
$ws.Range("L7").Value = 4333
$ws.Range("M7").Value = -14686.091
$ws.Range("N7").Value = -4557
$ws.Range("H40").Value = 3659.8
$ws.Range("I40").Value = 3659.8
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3659.8
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3523.8
$ws.Range("N40").Value = ""
$ws.Range("H55").Value = 263.22223
$ws.Range("I55").Value = 243.75
$ws.Range("J55").Value = 278.8
$ws.Range("K55").Value = 243.75
$ws.Range("L55").Value = 278.8
$ws.Range("M55").Value = -70.75
$ws.Range("N55").Value = -624.8
$ws.Range("H74").Value = 34188
$ws.Range("I74").Value = 35000
$ws.Range("J74").Value = 33646.668
$ws.Range("K74").Value = 35000
$ws.Range("L74").Value = 33646.668
$ws.Range("M74").Value = -34002
$ws.Range("N74").Value = -35642.668
$ws.Range("H77").Value = 34188
$ws.Range("I77").Value = 35000
$ws.Range("J77").Value = 33646.668
$ws.Range("K77").Value = 105000
$ws.Range("L77").Value = 100940.004
$ws.Range("M77").Value = -100008
$ws.Range("N77").Value = -110924.004
$ws.Range("H126").Value = 13926
$ws.Range("I126").Value = 14798.091
$ws.Range("J126").Value = 4333
$ws.Range("K126").Value = 44394.273
$ws.Range("L126").Value = 12999
$ws.Range("M126").Value = -41924.273
$ws.Range("N126").Value = -17939

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = ""
$ws.Range("H107").Value = 1671.9474
$ws.Range("I107").Value = 1759.8
$ws.Range("K107").Value = 5279.4
$ws.Range("M107").Value = -3359.4
